$d = $word.ActiveDocument

# The sentence originally reads (all one run):
#   " Century). Dennis (2023) reiterated how Southern Spain shifted from a
#   Periphery into a Core. The main reasons for this shift was due to"
#
# The edit keeps the "...into a Core" lead-in untouched and swaps the
# trailing clause for new wording:
#   ". The main reasons for this shift was due to"
#     -> " due to certain events. This event "

$old = ". The main reasons for this shift was due to"
$new = " due to certain events. This event "

$r = $d.Content
$found = $r.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

if (-not $found) {
    throw "Could not find the target sentence tail to replace."
}

# Word splits a run at the point of an edit, so the untouched "...into a
# Core" text and the freshly typed replacement end up as sibling <w:r>
# elements that share identical run formatting. Re-locate the text we just
# inserted and nudge a character property on (and back off) it so it is
# materialized as its own run rather than being folded back into the
# preceding run.
$r2 = $d.Content
$splitFound = $r2.Find.Execute($new, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $splitFound) {
    throw "Could not re-locate the inserted text to split it into its own run."
}
$r2.Bold = 1
$r2.Bold = 0
